$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the series. Insert a row above the
# current row 357 (shifting the existing rows 357:429 down to 358:430,
# which also grows the sheet's used range to row 430) and populate the
# freshly inserted row with the new record.
$ws.Rows("357:357").Insert()

$ws.Range("A357").Value = 10
$ws.Range("B357").Value = "Vega Modelo de Temuco"
$ws.Range("C357").Value = "La Araucanía"
$ws.Range("D357").Value = 45209
$ws.Range("E357").Value = 9
$ws.Range("F357").Value = 100112039
$ws.Range("G357").Value = "Ciboulette"
$ws.Range("H357").Value = "Sin especificar"
$ws.Range("I357").Value = "Primera"
$ws.Range("J357").Value = 30
$ws.Range("K357").Value = 7000
$ws.Range("L357").Value = 7000
$ws.Range("M357").Value = 7000
$ws.Range("N357").Value = "$/docena de atados"
$ws.Range("O357").Value = "Provincia de Cautín"
$ws.Range("P357").Value = 2333
$ws.Range("Q357").Value = 3
$ws.Range("R357").Value = "Hortaliza"
